$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clear the old "5 day moving" / label cells that are being replaced by
#    the new Return % / Correlation layout.
# ---------------------------------------------------------------------------
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("H11").ClearContents()

# ---------------------------------------------------------------------------
# 2. Row 1 headers
# ---------------------------------------------------------------------------
$ws.Range("F1").Value2 = "Return %"
$ws.Range("G1").Value2 = "Date"
$ws.Range("H1").Value2 = "CSI300"
$ws.Range("I1").Value2 = "ZhongZheng 500"
$ws.Range("J1").Value2 = "Chuang Ye Ban"
$ws.Range("L1").Value2 = "Correlation"
$ws.Range("M1").Value2 = "MV 5"

# ---------------------------------------------------------------------------
# 3. Column G - mirror of the Date column (A) for rows 2:16
# ---------------------------------------------------------------------------
$ws.Range("G2:G16").Formula = "=A2"
$ws.Range("G2:G16").Value2 = $ws.Range("G2:G16").Value2

# ---------------------------------------------------------------------------
# 4. Columns H / I - daily % return of CSI300 / ZhongZheng 500, rows 3:16
# ---------------------------------------------------------------------------
$ws.Range("H3").Formula = "=(B3/B2-1)*100"
$ws.Range("I3").Formula = "=(C3/C2-1)*100"
$ws.Range("H4:H16").Formula = "=(B4/B3-1)*100"
$ws.Range("I4:I16").Formula = "=(C4/C3-1)*100"

# Column J - daily % return of Chuang Ye Ban, rows 9:16 (filled right from I)
$ws.Range("J9:J16").Formula = "=(D9/D8-1)*100"

# ---------------------------------------------------------------------------
# 5. Column L - rolling correlation between H and I (5-day window), rows 7:16
# ---------------------------------------------------------------------------
$ws.Range("L7").Formula = "=CORREL(H3:H7,I3:I7)"
$ws.Range("L8:L16").Formula = "=CORREL(H4:H8,I4:I8)"

# Label under L6 identifying what L's correlation series represents
$ws.Range("L6").Value2 = "CSI300-ZhongZheng 500"

# ---------------------------------------------------------------------------
# 6. Columns M / N - rolling correlation between I&J and H&J, rows 13:16
# ---------------------------------------------------------------------------
$ws.Range("M13").Formula = "=CORREL(I9:I13,J9:J13)"
$ws.Range("M14").Formula = "=CORREL(I10:I14,J10:J14)"
$ws.Range("M15").Formula = "=CORREL(I11:I15,J11:J15)"
$ws.Range("M16").Formula = "=CORREL(I12:I16,J12:J16)"

$ws.Range("N13").Formula = "=CORREL(H9:H13,J9:J13)"
$ws.Range("N14:N16").Formula = "=CORREL(H10:H14,J10:J14)"

# Labels identifying the M / N correlation series
$ws.Range("M12").Value2 = "ZhongZheng 500-Chuang Ye Ban"
$ws.Range("N12").Value2 = "CSI300-Chuang Ye Ban"

Write-Host "values and formulas written"
